$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the (encrypted/obfuscated) text content of B14 to the new value.
$ws.Range("B14").Value = "U2FsdGVkX18/i8S14xFv/odbT8CTiiDUYFdJCnQ1YIpPHKsHgFx9UsHkbGmtDI7aeETe3bhfq+zQlQ3oQm0zDoBKlkGsRyVdGhTyw8/UYDNTp/mTJxpcwTawDhYhUKQ+1mmiZ9qmocMabuAcgh7IXcfyuWQZCCcc9dgDzNC5XSz3SRbCYPlT6BXWGfGUCQKd5w3ChWeRs723vrdzXvTAAzmY80iJ3Z1uY5Zo1rlFaAga/v2vzpzw+hYDOPP3dDGBObuqPfox7pbBMaBb9iBHupU6NxXydKY0y44aT3pJ+GT4c6AXvD8RYVygr1T65UDgeoE/kRk6N4ow5EyOv25E42BHD51hLF19OUAp1n6EqREKNKc2CMln3RfA42/VaCDv08I1WuHM/uoTCY90xnAw8g6esa6r/988qvMSAhj+8FZHT4ur3cWejaWxn08Fz5rBD3aSsO88Dd96Q8VoXAudLgBwlMsvFM4kM9gAQNg8sJY2k7B6fECzM4WDTu/syibmvoeq2MOY94UVl2bf6jAREVelg6wdwfBlFV9OHATflAJA+pkwd1Emy9enufs0kpnWsGybnjlIF6hfOg9fTISUWJR7282fbzms3hGEbXrcLJdA3diacyes+McI38+f1+2jSgycVESEGPlerVWeweSHOgWKzU7RT4bTMc9A0CmI9/ct4E+pweQTFPMqrVz6y+GSv9ePn6jki4v1vfJ/6YUDzBxsiEoKTfFxPm3HAegUeb2vLEuB2RmMd3QN7kImKopmLRnThNCDRRdp2Si0bJ1Z6/KoRPWa4rEbb/lkc+aH6GI="

# Move the active selection from B14 to B18 on the frozen pane, matching
# where the user was working after the edit.
$ws.Range("B18").Select()
